# Weekly update for the Hortaliza / Femacal de La Calera - Albahaca sheet.
# The underlying data block (rows 2-113) is a rolling window of weekly price
# records. This edit inserts two new weekly records (pushing the older rows
# down by one row each) and renumbers/reflows the table accordingly, exactly
# as reflected in the XML diff (dimension grows from R113 to R115).
#
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across every data row in this
# block (same market/category/unit/etc.), so only D (date), J (volume),
# K/L/M (min/max/weighted price) and P (price per kg) vary and need to be
# supplied explicitly for the two newly inserted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$staticCols = @(1,2,3,5,6,7,8,9,14,15,17,18)  # A,B,C,E,F,G,H,I,N,O,Q,R

# NOTE: named parameters do not bind reliably in this PowerShell-style
# runtime, so the helper below takes plain positional arguments.
function Insert-WeeklyRow($RowIndex, $Date, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {

    $ws.Rows.Item($RowIndex).Insert()

    $sourceRow = $RowIndex + 1
    foreach ($col in $staticCols) {
        $ws.Cells.Item($RowIndex, $col).Value = $ws.Cells.Item($sourceRow, $col).Value2
    }

    $ws.Cells.Item($RowIndex, 4).Value = $Date
    $ws.Cells.Item($RowIndex, 4).NumberFormat = $ws.Cells.Item($sourceRow, 4).NumberFormat

    $ws.Cells.Item($RowIndex, 10).Value = $Volumen
    $ws.Cells.Item($RowIndex, 11).Value = $PrecioMin
    $ws.Cells.Item($RowIndex, 12).Value = $PrecioMax
    $ws.Cells.Item($RowIndex, 13).Value = $PrecioProm
    $ws.Cells.Item($RowIndex, 16).Value = $PrecioKg
}

# New record inserted before the old row 70 (2022-01-06).
Insert-WeeklyRow 70 44567 130 4500 5000 4769 795

# New record inserted before what is now row 111 (old row 110, after the
# first insert shifted everything below row 70 down by one) (2022-01-07).
Insert-WeeklyRow 111 44568 130 4500 5000 4769 795
